# Weekly fruit/vegetable price update: permute the (Fecha, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
# tuple across the data rows (2..27) of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# mapping[targetRow] = sourceRow -- i.e. after the edit, row `targetRow`
# holds the values that row `sourceRow` had before the edit.
$mapping = @{
    2  = 4
    3  = 15
    4  = 10
    5  = 27
    6  = 3
    7  = 18
    8  = 20
    9  = 25
    10 = 11
    11 = 14
    12 = 12
    13 = 24
    14 = 22
    15 = 5
    16 = 8
    17 = 19
    18 = 16
    19 = 9
    20 = 26
    21 = 23
    22 = 6
    23 = 13
    24 = 2
    25 = 21
    26 = 17
    27 = 7
}

# Columns whose values move together as a row-tuple.
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# Snapshot all current values first, since the permutation reads from
# rows that will also be written.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowVals
}

foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($row, $col).Value2 = $srcVals[$col]
    }
}
